$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "c_E2g_E/c_E"
$ws.Range("B1").Value = "c_E2g_E/g_E"
$ws.Range("C1").Value = "c_H2g_H/c_H"
$ws.Range("D1").Value = "c_H2g_H/g_H"
